# Moved PostProcessing from Unity into main.sh
# Adds a "Power Strip" line item (under a new "ADDITIONS: made on 7/5/2017"
# section) to the budget sheet, right before the existing Total row, and
# extends the trailing blank rows by two so the sheet still ends with the
# same amount of blank padding below the Total.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Make room: push the old "Total" row (17) and everything below it
#    down by two rows so we can insert the new section header row and
#    the new line-item row in front of it.
# ------------------------------------------------------------------
$ws.Rows("17:18").Insert()

# ------------------------------------------------------------------
# 2. New row 18: the "Power Strip" line item (mirrors the other item
#    rows, e.g. row 12/13/... with style 5 on A, plain text on B,
#    quantity in C, currency in D/E and the link in F).
#    NOTE: values are written in the same order that the shared string
#    table gains them in the target workbook (Power Strip, then the
#    full name, then the link, then the section header) so the new
#    shared-string indices line up with the diff (45,46,47,48).
# ------------------------------------------------------------------
$ws.Range("A18").Value2 = "Power Strip"
$ws.Range("B18").Value2 = "6-Outlet Surge Protector Power Strip"
$ws.Range("C18").Value2 = 1
$ws.Range("D18").Value2 = 10.14
$ws.Range("E18").Formula = "=C18*D18"
$ws.Range("F18").Value2 = "https://www.amazon.com/AmazonBasics-6-Outlet-Surge-Protector-Power/dp/B00TP1C51M/ref=sr_1_5?ie=UTF8&qid=1499273840&sr=8-5&keywords=power+strip"

# ------------------------------------------------------------------
# 3. New row 17: the "ADDITIONS: made on 7/5/2017" section header,
#    styled the same bold way as the other ADDITIONS rows (8 and 11).
# ------------------------------------------------------------------
$ws.Range("A17").Value2 = "ADDITIONS: made on 7/5/2017"
$ws.Range("A17").Font.Bold = $true

# ------------------------------------------------------------------
# 4. Fix up the Total row (now row 19) so it sums through the new
#    item row. (The single 2-row insert above already pushed the
#    trailing blank rows down so the sheet now correctly ends at
#    row 35 / A1:F35, matching the target workbook.)
# ------------------------------------------------------------------
$ws.Range("E19").Formula = "=SUM(E2:E18)"

# ------------------------------------------------------------------
# 5. Restore the active cell/selection like the target workbook.
# ------------------------------------------------------------------
$ws.Range("B25").Select()
